# Update results table in "p_mw.xlsx" (res_bus active power results)
# per the recorded diff: column B (p_mw) values change and column K
# (shift-related column, values go to 0) for data rows 2..25, plus a
# few column I values shift by floating point rounding.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (bus 0)
$ws.Range("B2").Value = -33.79776967907753
$ws.Range("K2").Value = 0

# Row 3 (bus 1)
$ws.Range("B3").Value = -16.57589304253133
$ws.Range("K3").Value = 0

# Row 4 (bus 2)
$ws.Range("B4").Value = -3.739546201286657
$ws.Range("K4").Value = 0

# Row 5 (bus 3)
$ws.Range("B5").Value = 2.015120619451579
$ws.Range("K5").Value = 0

# Row 6 (bus 4)
$ws.Range("B6").Value = 4.127387506939471
$ws.Range("I6").Value = -36.68639292369173
$ws.Range("K6").Value = 0

# Row 7 (bus 5)
$ws.Range("B7").Value = 4.127387506939471
$ws.Range("I7").Value = -36.68639292369173
$ws.Range("K7").Value = 0

# Row 8 (bus 6)
$ws.Range("B8").Value = 3.178874857707797
$ws.Range("K8").Value = 0

# Row 9 (bus 7)
$ws.Range("B9").Value = -8.126929238241303
$ws.Range("I9").Value = -36.17910782839056
$ws.Range("K9").Value = 0

# Row 10 (bus 8)
$ws.Range("B10").Value = -34.49879428346088
$ws.Range("K10").Value = 0

# Row 11 (bus 9)
$ws.Range("B11").Value = -60.01034101119308
$ws.Range("K11").Value = 0

# Row 12 (bus 10)
$ws.Range("B12").Value = -71.76806989927343
$ws.Range("K12").Value = 0

# Row 13 (bus 11)
$ws.Range("B13").Value = -75.05619390936548
$ws.Range("I13").Value = -37.23736614058325
$ws.Range("K13").Value = 0

# Row 14 (bus 12)
$ws.Range("B14").Value = -71.24547154173379
$ws.Range("K14").Value = 0

# Row 15 (bus 13)
$ws.Range("B15").Value = -67.14130586339917
$ws.Range("I15").Value = -36.00561255416113
$ws.Range("K15").Value = 0

# Row 16 (bus 14)
$ws.Range("B16").Value = -65.6038481443392
$ws.Range("I16").Value = -36.72505219675369
$ws.Range("K16").Value = 0

# Row 17 (bus 15)
$ws.Range("B17").Value = -64.21287019300769
$ws.Range("I17").Value = -37.65916081090182
$ws.Range("K17").Value = 0

# Row 18 (bus 16)
$ws.Range("B18").Value = -57.02513233280061
$ws.Range("I18").Value = -37.56895584042385
$ws.Range("K18").Value = 0

# Row 19 (bus 17)
$ws.Range("B19").Value = -53.83633321965641
$ws.Range("I19").Value = -37.98697887434611
$ws.Range("K19").Value = 0

# Row 20 (bus 18)
$ws.Range("B20").Value = -61.47258285203563
$ws.Range("I20").Value = -37.9049457827268
$ws.Range("K20").Value = 0

# Row 21 (bus 19)
$ws.Range("B21").Value = -80.53185158283752
$ws.Range("K21").Value = 0

# Row 22 (bus 20)
$ws.Range("B22").Value = -92.55200030213092
$ws.Range("K22").Value = 0

# Row 23 (bus 21)
$ws.Range("B23").Value = -101.0304673417772
$ws.Range("K23").Value = 0

# Row 24 (bus 22)
$ws.Range("B24").Value = -103.8239395132402
$ws.Range("K24").Value = 0

# Row 25 (bus 23)
$ws.Range("B25").Value = -84.64377263253711
$ws.Range("I25").Value = -42
$ws.Range("K25").Value = 0

Write-Host "Updated p_mw results: columns B, K (and minor I roundings) for rows 2-25"
